# Applies the 6.4.1.2 "Water loss during transportation" update:
#  - Fix the title strings (remove stray period after "6.4.1.2")
#  - Update the 2022 (column P) figures for Kyrgyz Republic / Osh oblast
#  - Update the active selection to S3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title row (row 1): Kyrgyz / Russian / English captions.
# Only the Russian (B1) and English (C1) captions change, from
# "6.4.1.2. ..." to "6.4.1.2 ..." (dot after "2" removed).
$ws.Range("B1").Value = "6.4.1.2 Потери воды при транспортировке"
$ws.Range("C1").Value = "6.4.1.2 Percentage of water loss during transportation"

# Updated 2022 figures (column P)
$ws.Range("P5").Value = 2388        # Кыргыз Республикасы - Million cube meters table
$ws.Range("P10").Value = 335.3      # Ош облусу - Million cube meters table
$ws.Range("P16").Value = 27.3       # Кыргыз Республикасы - Percentage table
$ws.Range("P21").Value = 24.3       # Ош облусу - Percentage table

# Update the saved selection/active cell
$ws.Range("S3").Select()
